# Refresh the cryptos list (Price / Volume(1h) columns) with the latest
# scraped values, matching the "Updated cryptos list ... with GitHub Actions"
# commit. Also PEPE and dogwifhat swapped rank (rows 43/44) so their
# Coin/Link/Price/Volume columns are swapped too.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text format on numeric-looking D-column cells so values stay text (matching original inlineStr type)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D51").NumberFormat = "@"

# Apply updated cell values per diff
$ws.Range("D2").Value = '66.506.64'
$ws.Range("E2").Value = '  +1.17%  '
$ws.Range("D3").Value = '3.188.86'
$ws.Range("E3").Value = '  -0.12%  '
$ws.Range("E4").Value = '  -0.04%  '
$ws.Range("D5").Value = '604.29'
$ws.Range("E5").Value = '  +1.74%  '
$ws.Range("D6").Value = '156.10'
$ws.Range("E6").Value = '  +3.49%  '
$ws.Range("D7").Value = '0.999'
$ws.Range("E7").Value = '  -0.03%  '
$ws.Range("D8").Value = '3.188.51'
$ws.Range("E8").Value = '  -0.14%  '
$ws.Range("E9").Value = '  +2.78%  '
$ws.Range("E10").Value = '  -1.30%  '
$ws.Range("E11").Value = '  -4.59%  '
$ws.Range("E12").Value = '  +0.04%  '
$ws.Range("E13").Value = '  -2.06%  '
$ws.Range("D14").Value = '38.98'
$ws.Range("E14").Value = '  +2.00%  '
$ws.Range("D15").Value = '3.711.55'
$ws.Range("E15").Value = '  -0.18%  '
$ws.Range("D16").Value = '66.526.98'
$ws.Range("E16").Value = '  +0.95%  '
$ws.Range("D17").Value = '7.41'
$ws.Range("E17").Value = '  +1.45%  '
$ws.Range("D18").Value = '3.184.82'
$ws.Range("E18").Value = '  -0.43%  '
$ws.Range("E19").Value = '  +0.12%  '
$ws.Range("D20").Value = '514.65'
$ws.Range("E20").Value = '  +0.61%  '
$ws.Range("D21").Value = '15.47'
$ws.Range("E21").Value = '  -3.04%  '
$ws.Range("D22").Value = '0.735'
$ws.Range("E22").Value = '  +0.28%  '
$ws.Range("D23").Value = '8.15'
$ws.Range("E23").Value = '  +2.50%  '
$ws.Range("D24").Value = '14.92'
$ws.Range("E24").Value = '  -2.28%  '
$ws.Range("D25").Value = '84.68'
$ws.Range("E25").Value = '  -0.97%  '
$ws.Range("E26").Value = '  -0.04%  '
$ws.Range("D27").Value = '3.01'
$ws.Range("E27").Value = '  +0.10%  '
$ws.Range("D28").Value = '9.22'
$ws.Range("E28").Value = '  -0.39%  '
$ws.Range("E29").Value = '  +6.78%  '
$ws.Range("D30").Value = '3.06'
$ws.Range("E30").Value = '  +6.61%  '
$ws.Range("D31").Value = '7.05'
$ws.Range("E31").Value = '  +5.94%  '
$ws.Range("D32").Value = '28.12'
$ws.Range("E32").Value = '  +0.12%  '
$ws.Range("E33").Value = '  -1.62%  '
$ws.Range("E34").Value = '  -0.10%  '
$ws.Range("D35").Value = '6.55'
$ws.Range("E35").Value = '  -0.85%  '
$ws.Range("D36").Value = '513.63'
$ws.Range("E36").Value = '  +6.66%  '
$ws.Range("D37").Value = '54.74'
$ws.Range("E37").Value = '  -1.07%  '
$ws.Range("D38").Value = '0.0896'
$ws.Range("E38").Value = '  -2.25%  '
$ws.Range("D39").Value = '0.0424'
$ws.Range("E39").Value = '  -0.16%  '
$ws.Range("E40").Value = '  +6.01%  '
$ws.Range("D41").Value = '8.86'
$ws.Range("E41").Value = '  -1.11%  '
$ws.Range("E42").Value = '  +4.65%  '
$ws.Range("B43").Value = 'dogwifhat'
$ws.Range("C43").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range("D43").Value = '2.86'
$ws.Range("E43").Value = '  -4.33%  '
$ws.Range("B44").Value = 'PEPE'
$ws.Range("C44").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D44").Value = '0.0₃0676'
$ws.Range("E44").Value = '  +6.74%  '
$ws.Range("E45").Value = '  -0.76%  '
$ws.Range("D46").Value = '2.859.70'
$ws.Range("E46").Value = '  -5.41%  '
$ws.Range("D47").Value = '28.42'
$ws.Range("E47").Value = '  -2.39%  '
$ws.Range("D48").Value = '2.41'
$ws.Range("E48").Value = '  +4.71%  '
$ws.Range("E49").Value = '  +0.00%  '
$ws.Range("E50").Value = '  +0.74%  '
$ws.Range("D51").Value = '2.60'
$ws.Range("E51").Value = '  +6.01%  '
